# ------------------------------------------------------------------
# covid-vulnerability-index / tabs / prep-code-stab5.xlsx
# "added zipcode-level NEVI, county --> borough, NVI --> NEVI"
#
# - rename Sheet1 -> tab_fig
# - add a new sheet "fig" after tab_fig, and make it the active sheet
# - populate "fig" with the subdomain / subdomain-label lookup table
#   (columns A-D plus the helper column F), each with its own formula
# - tidy up the selection left behind on tab_fig
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$q  = [char]34   # a literal double-quote, used inside the R-code formula strings

# --- rename the original sheet ------------------------------------------------
$tabFig = $wb.Worksheets.Item(1)
$tabFig.Name = "tab_fig"

# --- add the new "fig" sheet right after tab_fig ------------------------------
$fig = $wb.Worksheets.Add($null, $tabFig)
$fig.Name = "fig"

# --- column widths (cosmetic, bestFit-like) -----------------------------------
$fig.Columns.Item(1).ColumnWidth = 38.499
$fig.Columns.Item(2).ColumnWidth = 13.499
$fig.Columns.Item(3).ColumnWidth = 18.054
$fig.Columns.Item(4).ColumnWidth = 18.054
$fig.Columns.Item(5).ColumnWidth = 18.054

# --- header row ----------------------------------------------------------------
$fig.Cells.Item(1,1).Value = "subdomain"
$fig.Cells.Item(1,2).Value = "subdomain label"
$fig.Cells.Item(1,3).Value = "subdomain label code"
$fig.Cells.Item(1,4).Value = "row_no"
$fig.Cells.Item(1,6).Value = "OLD"

# --- lookup table: subdomain / label / row number -----------------------------
# (A = subdomain *_median_iqr name, B = display label, D = row number)
$lookup = @(
    @("score_demo_age_median_iqr",              "Age"),
    @("score_demo_femaleled_median_iqr",         "Female-Led Households"),
    @("score_demo_immigration_median_iqr",       "Immigration"),
    @("score_demo_disability_median_iqr",        "Disability"),
    @("score_demo_singleparent_median_iqr",      "Single Parent Households"),
    @("score_demo_mobility_median_iqr",          "Mobility"),
    @("score_demo_livealone_median_iqr",         "Social Isolation"),
    @("score_economic_incomepoverty_median_iqr", "Income and Poverty"),
    @("score_economic_servicemanual_median_iqr", "Occupation"),
    @("score_economic_gini_median_iqr",          "Income Inequality"),
    @("score_economic_employment_median_iqr",    "Unemployment"),
    @("score_economic_education_median_iqr",     "Education"),
    @("score_economic_vehicleavail_median_iqr",  "Vehicle Availability"),
    @("score_residential_popdensity_median_iqr",    "Population Density"),
    @("score_residential_groupquarters_median_iqr", "Group Quarters"),
    @("score_residential_occperroom_median_iqr",    "Occupants Per Room"),
    @("score_residential_structage_median_iqr",     "Age of Housing Structure"),
    @("score_residential_structattach_median_iqr",  "Units in Housing Structure"),
    @("score_residential_move1yr_median_iqr",       "Changing Residence"),
    @("score_residential_vacancy_median_iqr",       "Vacancy"),
    @("score_healthstatus_lifestyle_median_iqr",      "Unhealthy Behaviors"),
    @("score_healthstatus_condition_median_iqr",      "Health Outcomes"),
    @("score_healthstatus_preventive_median_iqr",     "Prevention Practices"),
    @("score_healthstatus_lackinsurance_median_iqr",  "Health Insurance Access")
)

for ($i = 0; $i -lt $lookup.Count; $i++) {
    $r = $i + 2
    $subdomain = $lookup[$i][0]
    $label     = $lookup[$i][1]
    $rowNo     = $i + 1

    $fig.Cells.Item($r,1).Value = $subdomain
    $fig.Cells.Item($r,2).Value = $label
    $fig.Cells.Item($r,4).Value = $rowNo

    # C: "'<label>',"
    $fig.Cells.Item($r,3).Formula = "=" + $q + "'" + $q + "&B$r&" + $q + "'," + $q

    # F: "subdomain == '<subdomain>' ~ subdomain_label_vector[<row_no>],"
    $fig.Cells.Item($r,6).Formula = "=" + $q + "subdomain == '" + $q + "&A$r&" + $q + "' ~ subdomain_label_vector[" + $q + "&D$r&" + $q + "]," + $q
}

# --- selection / active-sheet bookkeeping --------------------------------------
$tabFig.Range("B17").Select() | Out-Null
$fig.Range("E7").Select() | Out-Null
$fig.Activate() | Out-Null
